$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet "page-1_table-1": row 2 date header, shifted from 22/07-28/07 to
# 05/08-11/08, now stored with an explicit Text number format.
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "05/08"
$ws1.Range("B2").Value = "06/08"
$ws1.Range("C2").Value = "07/08"
$ws1.Range("D2").Value = "08/08"
$ws1.Range("E2").Value = "09/08"
$ws1.Range("F2").Value = "10/08"
$ws1.Range("G2").Value = "11/08"
$ws1.Range("A2:G2").NumberFormat = "@"
$ws1.Range("H2").NumberFormat = "@"
$ws1.Range("M11").Select()

# ---------------------------------------------------------------------------
# Sheet "page-1_table-2": sample employee schedule row filled in, and the
# "Slot N" helper labels moved up so they read consecutively on rows 2-8.
# ---------------------------------------------------------------------------
$ws2.Range("B1").Value = "SE1824-IoT102t at P.005"
$ws2.Range("B2").Value = "SE1824-IoT102t at P.005"
$ws2.Range("B3").Value = "SE1824-IoT102t at P.005"
$ws2.Range("B4").Value = ""

$ws2.Range("A2").Value = "Slot 2"
$ws2.Range("A3").Value = "Slot 3"
$ws2.Range("A4").Value = "Slot 4"
$ws2.Range("A5").Value = "Slot 5"
$ws2.Range("A6").Value = "Slot 6"
$ws2.Range("A7").Value = "Slot 7"
$ws2.Range("A8").Value = "Slot 8"
$ws2.Range("A9").Value = ""
$ws2.Range("A11").Value = ""
$ws2.Range("A13").Value = ""
$ws2.Range("A14").Value = ""

# E2 loses the highlighted (filled) look it had before, matching plain A-col style
$ws2.Range("A1").Copy()
$ws2.Range("E2").PasteSpecial(-4122)

# A5 becomes the new "currently highlighted" slot row (style previously on A9/C9)
$ws2.Range("D9").Copy()
$ws2.Range("A5").PasteSpecial(-4122)

# A9 and C9 drop back to the plain (unhighlighted) style
$ws2.Range("B9").Copy()
$ws2.Range("A9").PasteSpecial(-4122)
$ws2.Range("C9").PasteSpecial(-4122)

$ws2.Application.CutCopyMode = $false

$ws2.Range("B7").Select()
$ws2.Application.ActiveWindow.Zoom = 85
